# grafikai javítások, alapértelmezett verseny az eredményfelvételnél
#
# Adds a "Kategoria" column (W), replaces the sample row-2 data with a new
# competitor, fills in row 3 (previously blank) with another competitor and
# appends a brand new row 4 with a third competitor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a value into a cell but force it to be stored as TEXT
# even when it "looks" numeric (e.g. "35.0", "25", "28"), mirroring how the
# source export always uses inline strings. We briefly mark the cell as
# Text, assign the value, then drop back to the Normal style so no stray
# number-format styling is left behind on the cell. ------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- New header column W1 ("Kategoria"), copying the look of the other
# header cells (bold/centered/bordered) from V1. ---------------------------
$ws.Range("V1").Copy($ws.Range("W1"))
$ws.Range("W1").Value = "Kategoria"

# --- Row 2: overwrite with the new competitor's data -----------------------
# Only cells that actually held a previous value need clearing explicitly
# (D2, E2, F2, J2, Q2, V2); the remaining cells in the row are already blank
# placeholders in the source sheet and are simply left untouched.
Set-TextValue $ws.Range("A2") "35.0"
$ws.Range("B2").Value = "Csefkó Pál"
$ws.Range("C2").Value = "Csikóvári Para Sport Egyesület"
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("H2").Value = "2025-09-13 20:47"
$ws.Range("J2").ClearContents()
$ws.Range("Q2").ClearContents()
$ws.Range("V2").ClearContents()

# --- Row 3: previously blank, now filled in with a second competitor ------
Set-TextValue $ws.Range("A3") "25"
$ws.Range("B3").Value = "Dr. Seffer István"
$ws.Range("C3").Value = "Szent Hubertus Szituációs Lövészklub"
$ws.Range("E3").Value = "o"
$ws.Range("H3").Value = "2025-09-14 22:14"
$ws.Range("V3").Value = "VID_00001"

# --- Row 4: brand new row for a third competitor ---------------------------
Set-TextValue $ws.Range("A4") "28"
$ws.Range("B4").Value = "Vivert János"
$ws.Range("C4").Value = "Szent Hubertus Szituációs Lövészklub"
$ws.Range("V4").Value = "VID_00001"
